$d = $word.ActiveDocument

function Remove-TrailingSpaceRun($headingText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $full = $p.Range.Text
        # Paragraph text includes the trailing paragraph mark (chr 13).
        if ($full.Length -ge 2 -and $full.Substring(0, $full.Length - 1) -eq ($headingText + " ")) {
            $r = $p.Range
            # Delete just the trailing space character, leaving the paragraph mark.
            $d.Range($r.End - 2, $r.End - 1).Delete()
            return
        }
    }
}

# --- Remove the lone trailing-space run after the "Water Skiing" heading ---
Remove-TrailingSpaceRun("Water Skiing")

# --- Remove the lone trailing-space run after the "Photography" heading ---
Remove-TrailingSpaceRun("Photography")

# --- Resize the Memorials table's two grid columns (values are in points;
#     the engine persists them as twentieths of a point / twips, so
#     270pt -> 5400 dxa and 126pt -> 2520 dxa) ---
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 270
$t.Columns.Item(2).Width = 126
